$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$urls = @(
    "https://www.nikkei.com/article/DGXMZO53049580W9A201C1DTA000",
    "https://www.nikkei.com/article/DGXZZO53052440W9A201C1000000",
    "https://www.nikkei.com/article/DGXMZO53052770W9A201C1916M00",
    "https://www.nikkei.com/article/DGXMZO53050320W9A201C1TJC000",
    "https://www.nikkei.com/article/DGXMZO53051890W9A201C1L83000",
    "https://www.nikkei.com/article/DGXMZO53049880W9A201C1000000",
    "https://www.nikkei.com/article/DGXMZO53049690W9A201C1TJC000",
    "https://www.nikkei.com/article/DGXLASFL06HTR_W9A201C1000000",
    "https://www.nikkei.com/article/DGXLMS1333H5FDW9A201C1000000",
    "https://www.nikkei.com/article/DGXLASFL06HRJ_W9A201C1000000",
    "https://www.nikkei.com/article/DGKKZO52998200V01C19A2BE0P00",
    "https://www.nikkei.com/article/DGXMZO53042440W9A201C1000000",
    "https://www.nikkei.com/article/DGXMZO53039790W9A201C1000000",
    "https://www.nikkei.com/article/DGXZZO53034250W9A201C1000000",
    "https://www.nikkei.com/article/DGXMZO53033960W9A201C1000000",
    "https://www.nikkei.com/article/DGXMZO52980070V01C19A2962M00",
    "https://www.nikkei.com/article/DGXMZO51480330Y9A021C1I00000",
    "https://www.nikkei.com/article/DGKKZO53023150V01C19A2FFE000",
    "https://www.nikkei.com/article/DGXMZO53016470V01C19A2FFE000",
    "https://www.nikkei.com/article/DGXMZO53015600V01C19A2TJ2000"
)

$row = 3
foreach ($url in $urls) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $url
    $cell.Style = "Normal"
    $row = $row + 1
}
